$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2193.52
$ws.Range("I28").Value = 2282.4736
$ws.Range("J28").Value = 1911.8334
$ws.Range("K28").Value = 2282.4736
$ws.Range("L28").Value = 1911.8334
$ws.Range("M28").Value = -1797.4736
$ws.Range("N28").Value = -2881.8334

$ws.Range("H39").Value = 1942
$ws.Range("I39").Value = 737.6
$ws.Range("K39").Value = 2212.8
$ws.Range("M39").Value = -1916.8

$ws.Range("H70").Value = 6772.522
$ws.Range("I70").Value = 2415.5881
$ws.Range("J70").Value = 19117.166
$ws.Range("K70").Value = 7246.7643
$ws.Range("L70").Value = 57351.49800000001
$ws.Range("M70").Value = -6976.7643
$ws.Range("N70").Value = -57891.49800000001

$ws.Range("H73").Value = 6772.522
$ws.Range("I73").Value = 2415.5881
$ws.Range("J73").Value = 19117.166
$ws.Range("K73").Value = 7246.7643
$ws.Range("L73").Value = 57351.49800000001
$ws.Range("M73").Value = -6310.7643
$ws.Range("N73").Value = -59223.49800000001

$ws.Range("H80").Value = 1097.8
$ws.Range("I80").Value = 496.33334
$ws.Range("K80").Value = 1489.00002
$ws.Range("M80").Value = -491.0000199999999

$ws.Range("H83").Value = 1097.8
$ws.Range("I83").Value = 496.33334
$ws.Range("K83").Value = 4467.00006
$ws.Range("M83").Value = 524.9999399999997

$ws.Range("H86").Value = 3699.3333
$ws.Range("I86").Value = 3699
$ws.Range("J86").Value = 3700
$ws.Range("K86").Value = 3699
$ws.Range("L86").Value = 3700
$ws.Range("M86").Value = -2576
$ws.Range("N86").Value = -5946

$ws.Range("H89").Value = 3699.3333
$ws.Range("I89").Value = 3699
$ws.Range("J89").Value = 3700
$ws.Range("K89").Value = 18495
$ws.Range("L89").Value = 18500
$ws.Range("M89").Value = -12879
$ws.Range("N89").Value = -29732

$ws.Range("H98").Value = 12184.111
$ws.Range("I98").Value = 930.4
$ws.Range("J98").Value = 26251.25
$ws.Range("K98").Value = 930.4
$ws.Range("L98").Value = 26251.25
$ws.Range("M98").Value = 567.6
$ws.Range("N98").Value = -29247.25

$ws.Range("H122").Value = 12184.111
$ws.Range("I122").Value = 930.4
$ws.Range("J122").Value = 26251.25
$ws.Range("K122").Value = 2791.2
$ws.Range("L122").Value = 78753.75
$ws.Range("M122").Value = -341.1999999999998
$ws.Range("N122").Value = -83653.75

$ws.Range("H134").Value = 35998.688
$ws.Range("J134").Value = 35998.688
$ws.Range("L134").Value = 35998.688
$ws.Range("N134").Value = -46138.688

$ws.Range("H137").Value = 1669.6222
$ws.Range("I137").Value = 1538.4286
$ws.Range("K137").Value = 4615.2858
$ws.Range("M137").Value = -2065.2858


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3049.3823
$ws.Range("I32").Value = 2215.3438
$ws.Range("K32").Value = 2215.3438
$ws.Range("M32").Value = -1928.3438

$ws.Range("H61").Value = 2091.5334
$ws.Range("I61").Value = 1788.7407
$ws.Range("K61").Value = 1788.7407
$ws.Range("M61").Value = -1576.7407

$ws.Range("H88").Value = 9086.923000000001
$ws.Range("I88").Value = 1459.75
$ws.Range("J88").Value = 12476.777
$ws.Range("K88").Value = 1459.75
$ws.Range("L88").Value = 12476.777
$ws.Range("M88").Value = -1053.75
$ws.Range("N88").Value = -13288.777

$ws.Range("H91").Value = 9086.923000000001
$ws.Range("I91").Value = 1459.75
$ws.Range("J91").Value = 12476.777
$ws.Range("K91").Value = 1459.75
$ws.Range("L91").Value = 12476.777
$ws.Range("M91").Value = -55.75
$ws.Range("N91").Value = -15284.777

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H136").Value = 2091.5334
$ws.Range("I136").Value = 1788.7407
$ws.Range("K136").Value = 5366.2221
$ws.Range("M136").Value = -2816.2221


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8589.956
$ws.Range("I86").Value = 3281.7222
$ws.Range("J86").Value = 27699.6
$ws.Range("K86").Value = 3281.7222
$ws.Range("L86").Value = 27699.6
$ws.Range("M86").Value = -2158.7222
$ws.Range("N86").Value = -29945.6

$ws.Range("H89").Value = 8589.956
$ws.Range("I89").Value = 3281.7222
$ws.Range("J89").Value = 27699.6
$ws.Range("K89").Value = 16408.611
$ws.Range("L89").Value = 138498
$ws.Range("M89").Value = -10792.611
$ws.Range("N89").Value = -149730

$ws.Range("H99").Value = 2346.6191
$ws.Range("I99").Value = 2120.5264
$ws.Range("J99").Value = 4494.5
$ws.Range("K99").Value = 2120.5264
$ws.Range("L99").Value = 4494.5
$ws.Range("M99").Value = -622.5264000000002
$ws.Range("N99").Value = -7490.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2313.625
$ws.Range("I31").Value = 1832.1111
$ws.Range("J31").Value = 3758.1667
$ws.Range("K31").Value = 1832.1111
$ws.Range("L31").Value = 3758.1667
$ws.Range("M31").Value = -1537.1111
$ws.Range("N31").Value = -4348.1667

$ws.Range("H34").Value = 2313.625
$ws.Range("I34").Value = 1832.1111
$ws.Range("J34").Value = 3758.1667
$ws.Range("K34").Value = 1832.1111
$ws.Range("L34").Value = 3758.1667
$ws.Range("M34").Value = -1630.1111
$ws.Range("N34").Value = -4162.1667

$ws.Range("H132").Value = 3658.28
$ws.Range("I132").Value = 2073.6
$ws.Range("K132").Value = 6220.799999999999
$ws.Range("M132").Value = -3690.799999999999


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 428.83334
$ws.Range("I5").Value = 430.7
$ws.Range("J5").Value = 419.5
$ws.Range("K5").Value = 1292.1
$ws.Range("L5").Value = 1258.5
$ws.Range("M5").Value = -1180.1
$ws.Range("N5").Value = -1482.5

$ws.Range("H135").Value = 428.83334
$ws.Range("I135").Value = 430.7
$ws.Range("J135").Value = 419.5
$ws.Range("K135").Value = 3876.3
$ws.Range("L135").Value = 3775.5
$ws.Range("M135").Value = -1341.3
$ws.Range("N135").Value = -8845.5


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2342.9546
$ws.Range("I132").Value = 2659.5
$ws.Range("K132").Value = 7978.5
$ws.Range("M132").Value = -5448.5


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2379.4722
$ws.Range("I132").Value = 2019.0667
$ws.Range("K132").Value = 6057.2001
$ws.Range("M132").Value = -3527.2001


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5984.8184
$ws.Range("I81").Value = 6233.8
$ws.Range("K81").Value = 12467.6
$ws.Range("M81").Value = -11406.6

$ws.Range("H84").Value = 5984.8184
$ws.Range("I84").Value = 6233.8
$ws.Range("K84").Value = 62338
$ws.Range("M84").Value = -57034

